$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-G (numeric + group member) for new rows 16-19
$ws.Cells.Item(16, 1).Value2 = 13
$ws.Cells.Item(16, 2).Value2 = 'Linthe'
$ws.Cells.Item(16, 3).Value2 = 11
$ws.Cells.Item(16, 4).Value2 = 53
$ws.Cells.Item(16, 5).Value2 = 5
$ws.Cells.Item(16, 6).Value2 = 78
$ws.Cells.Item(16, 7).Value2 = 25
$ws.Cells.Item(17, 1).Value2 = 15
$ws.Cells.Item(17, 2).Value2 = 'Linthe'
$ws.Cells.Item(17, 3).Value2 = 3
$ws.Cells.Item(17, 4).Value2 = 41
$ws.Cells.Item(17, 5).Value2 = 4
$ws.Cells.Item(17, 6).Value2 = 61
$ws.Cells.Item(17, 7).Value2 = 22
$ws.Cells.Item(18, 1).Value2 = 15
$ws.Cells.Item(18, 2).Value2 = 'Linthe'
$ws.Cells.Item(18, 3).Value2 = 4
$ws.Cells.Item(18, 4).Value2 = 46
$ws.Cells.Item(18, 5).Value2 = 1
$ws.Cells.Item(18, 6).Value2 = 51
$ws.Cells.Item(18, 7).Value2 = 22
$ws.Cells.Item(19, 1).Value2 = 14
$ws.Cells.Item(19, 2).Value2 = 'Linthe'
$ws.Cells.Item(19, 3).Value2 = 1
$ws.Cells.Item(19, 4).Value2 = 63
$ws.Cells.Item(19, 5).Value2 = 7
$ws.Cells.Item(19, 6).Value2 = 98
$ws.Cells.Item(19, 7).Value2 = 24

# Column H (Chat history) for new rows 16-19
$ws.Cells.Item(16, 8).Value2 = '[''Linthe Van Rooij: A woman dies after lighting a lighter, and doesnâ€™t even need a funeral.'', ''Thijs : Did she die in a fire or die from the smoke of the fire'', ''Linthe Van Rooij: No'', ''Thijs : Was she already burried?'', ''Linthe Van Rooij: Yes'', ''Thijs : Was she already presumed to be dead?'', ''Linthe Van Rooij: No'', ''Thijs : Did she fall in a coffin?'', ''Linthe Van Rooij: No'', ''Thijs : But she is in a coffin?'', ''Linthe Van Rooij: Yes'', ''Thijs : She lit the lighter to see that she was buried?'', ''Linthe Van Rooij: No'', ''Thijs : Was she mistaken for someone else?'', ''Linthe Van Rooij: No'', ''Thijs : Was she in the coffin by accident?'', ''Linthe Van Rooij: No'', ''Thijs : Was she looking for someone in the coffin?'', ''Linthe Van Rooij: No'', ''Thijs : Did people already think that she was dead?'', ''Linthe Van Rooij: No'', ''Thijs : Did she fall asleep in the coffin?'', ''Linthe Van Rooij: No'', ''Thijs : Was there anyone else involved in her death?'', ''Linthe Van Rooij: Yes'', ''Thijs : Did that person try to kill/murder her or accidentally kill her?'', ''Linthe Van Rooij: No'', ''Thijs : Was that person alive?'', ''Linthe Van Rooij: No'', ''Thijs : Did she suffocate?'', ''Linthe Van Rooij: Yes'', ''Thijs : Did she dig the hole that she was in herself?'', ''Linthe Van Rooij: No'', ''Thijs : Hint'', ''Linthe Van Rooij: Ask about the occupation of the other person involved'', ''Thijs : Was the other person a gravedigger or involved at the graveyard?'', ''Linthe Van Rooij: Yes'', ''Thijs : Was there anyone else involved'', ''Linthe Van Rooij: No'', ''Thijs : Was the gravedigger digging her grave?'', ''Linthe Van Rooij: No'', ''Thijs : Was the gravedigger already dead when she met him?'', ''Linthe Van Rooij: No'', ''Thijs : Did they fight and she killed him?'', ''Linthe Van Rooij: No'', ''Thijs : Did she know him?'', ''Linthe Van Rooij: Yes'', ''Thijs : Was she related to him?'', ''Linthe Van Rooij: No'', ''Thijs : Were they partners'', ''Linthe Van Rooij: No'', ''Linthe Van Rooij: Why did she lit the lighter inside the coffin?'', ''Thijs : To set something on fire or to signal something?'', ''Linthe Van Rooij: No'', ''Thijs : To light a cigarette or to suffocate herself?'', ''Linthe Van Rooij: No'', ''Thijs : Was the other person with her inside or the coffin?'', ''Linthe Van Rooij: Yes'', ''Thijs : Did she light it to kill him or to do anything related to him?'', ''Linthe Van Rooij: No'', ''Thijs : Did she take a shortcut home through the graveyard?'', ''Linthe Van Rooij: No'', ''Thijs : Did the guy fall in the coffin?'', ''Linthe Van Rooij: No'', ''Linthe Van Rooij: Ask about where she was before'', ''Thijs : Was she at a funeral?'', ''Linthe Van Rooij: No'', ''Thijs : At a social event?'', ''Linthe Van Rooij: No'', ''Thijs : Was she visiting someone at the graveyard or was she running?'', ''Linthe Van Rooij: No'', ''Thijs : Was she walking'', ''Linthe Van Rooij: No'', ''Thijs : Was she scared of the gravedigger'', ''Linthe Van Rooij: No'', ''Thijs : Was she at home'', ''Linthe Van Rooij: No'', ''Thijs : Anything sport related or food related'', ''Linthe Van Rooij: No'', ''Thijs : Was she drunk or on drugs'', ''Linthe Van Rooij: No'', ''Thijs : Was she shopping or in town'', ''Linthe Van Rooij: No'', ''Thijs : Was it night'', ''Linthe Van Rooij: Not relevant'', ''Thijs : Was she driving around'', ''Linthe Van Rooij: No'', ''Thijs : Hint'', ''Linthe Van Rooij: She was a criminal/convict'', ''Thijs : Did she just rob the house of the gravedigger'', ''Linthe Van Rooij: No'', ''Thijs : Did she rob someone who died'', ''Linthe Van Rooij: No'', ''Thijs : Did she just rob someone'', ''Linthe Van Rooij: No'', ''Thijs : Did she just murder someone?'', ''Linthe Van Rooij: No'', ''Thijs : Was she escaping prison'', ''Linthe Van Rooij: Yes'', ''Thijs : And she was hiding in the grave to escape the police?'', ''Linthe Van Rooij: Yes'', ''Thijs : And the gravedigger was an accomplice?'', ''Linthe Van Rooij: Yes'', ''Thijs : So they tried to hide from the police by hiding in the grave but they got locked in'', ''Linthe Van Rooij: No'', ''Thijs : She tried to hide from the police in the grave that the gravedigger made for her, but he was also dead in the coffin so they were trapped, she lit the lighter to confirm that he was dead'', ''Linthe Van Rooij: She expected another person to be in the grave next to her'', ''Thijs : She tried to escape the prison through the coffin of someone else, and the gravedigger was supposed to get her out, but the coffin was for the gravedigger thus he could not get her out, she used the lighter to see the face of the person that was in the coffin with her'', ''Linthe Van Rooij: Allrightyy'']'
$ws.Cells.Item(17, 8).Value2 = '[''Linthe Van Rooij: A man goes into a cafe and orders a glass of coke. The woman behind the bar grabs a knife and pointed it towards the man. He thanks her and leaves.'', ''Lynn van Zeijts: is one of them going to die after she leaves?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: do they have an important conversation after he received the coke?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: does the men also has a weapon?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: will he take the knife of the women?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: is she leaving with the knife in her hand?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: did the men ask the women for a knife?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: is one the men or the women getting injured?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: are there more people involved except the men & women?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: is it important to know why the women decides to leave?'', ''Linthe Van Rooij: Yes'', ''Linthe Van Rooij: False assumption woman leaves'', ''Lynn van Zeijts: so HE leaves? (not she)'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: did the  men pay for his coke?'', ''Lynn van Zeijts: is the relation between the men & women important?'', ''Linthe Van Rooij: Not relevant'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: did the women aimed to threaten him with the knife?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: did the women ask him to do something with the knife because she did not want to do something?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: did the men insult the women?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: did the men do something which can be seen considerd as ethically wrong?'', ''Linthe Van Rooij: No'', ''Linthe Van Rooij: Ask about why he orders something at the bar'', ''Lynn van Zeijts: did he order a coke to celebrate something?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: did he order something because had some medical emergency?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: did he order something for a negative reason?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: was this negative reason related to the women?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: was this related to his work?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: to his family live?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: is the fact that he order coke (and not another drink) important to understand the negative event?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: did he revelaed a secret?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: did the men receive his coke before she grabbed the knife?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: did the men threatend the women as well?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: so did the women refused to give him the coke?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: did the negative event happended shortly before he entered the cafe?'', ''Linthe Van Rooij: No, it is still happening'', ''Lynn van Zeijts: did they got into a fight?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: did the men drink anything before he order the coke?'', ''Linthe Van Rooij: No'', ''Linthe Van Rooij: The drink could serve as a remedy for the event'', ''Lynn van Zeijts: were there any other elements or objects involved, caused the negative event?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: is the coke dangerous for him?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: is the negative event over once he leaves the cafe?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: is the women trying to help the men by showing the knife?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: does the men wants to kill himself?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: does the women wants to kill or injure him?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: does the women wants to break or open anything related to the men (rope, hair, a bag ect)?'', ''Linthe Van Rooij: No'', ''Linthe Van Rooij: She only wants to scare him'', ''Lynn van Zeijts: to protect herself?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: to protect and help him?'', ''Linthe Van Rooij: Ask these questions seperately'', ''Lynn van Zeijts: is the women trying to help the men by scaring him with the knife?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: Did the men had a hiccup, and therefore asked for a coke. However, this did not work. Therefore the women tried to scared him with a knife. As a result, the hiccup went away and the men thanked her for her help.'', ''Linthe Van Rooij: Yes!!!!!!! You got it'']'
$ws.Cells.Item(18, 8).Value2 = '[''Linthe Van Rooij: A man climbs out of a large vehicle and drowned himself'', ''Lynn van Zeijts: is the vehicle still on land when he climbs out the vehicle?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: is the vehicle close to water as well?'', ''Linthe Van Rooij: Not relevant'', ''Lynn van Zeijts: were there any other people in the vehicle as well?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: was this vehicle an airplane?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: is this vehicle used for public transport?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: is it relevant to know the type of vehicle?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: is the vehicle made to be on land?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: is it a tractor?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: did he drown himself on purpose?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: was he a farmer?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: did something negative happened to his farming practices?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: did something negative happened at all?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: does this has anything to do with the tractor?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: was the tractor broken?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: was the land, in which the tractor was located at that moment boggy?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: did the reason that he climbed out of the tractor had anything to do with the tractor?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: was the water in which he drowned himself close to the tractor?'', ''Linthe Van Rooij: Not relevant'', ''Lynn van Zeijts: were there any other objects or people involved in the negative event?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: were there any people in the tractor?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: outside the tractor?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: did he wanted to go to the other people?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: did the drowning happened right after he climbed out of the tractor?'', ''Linthe Van Rooij: Not relevant'', ''Lynn van Zeijts: did the other people influenced him in stepping out of the tractor?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: was it dangerous for him to stay in the vehicle?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: were the other people also with tractors? n'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: is it important to know how he exactly climbed out of the tractor?'', ''Linthe Van Rooij: No'', ''Linthe Van Rooij: Ask about the relationship with the people involved'', ''Lynn van Zeijts: did the other people wanted him to drown himself?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: did the other people wanted to save him?'', ''Lynn van Zeijts: did he know the other people before?'', ''Linthe Van Rooij: Yes'', ''Linthe Van Rooij: Not relevant'', ''Lynn van Zeijts: was it an accident that he drowned?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: did the other people saw something which the men could not see (when he was still in the tractor)?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: were the other people also farmers?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: were the other people angry at him?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: did the know eachother for a long time?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: were they employees?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: family?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: were they close to the tractor?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: was there a fight or any other negative event between the men and his family?'', ''Linthe Van Rooij: No'', ''Lynn van Zeijts: were there any childern involved?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: were the childeren outside the tractor?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: did he climbed out the tractor for his childern?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: was his kid in a dangerous situation?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: did he aimed to save his kid by climbing out of the tractor?'', ''Lynn van Zeijts: did he cause the dangerous situation?'', ''Linthe Van Rooij: No'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: did he accidently killed his kid with his tractor?'', ''Linthe Van Rooij: Yes'', ''Lynn van Zeijts: so he accidently he run over the kid with his tractor because he did not saw him because they were playing the field so he drowned himself somewhere afterwards because it made him so depressed'']'
$ws.Cells.Item(19, 8).Value2 = '[''Linthe Van Rooij: A man and a woman are seated at a table. On the table between them is a game of checkers and a pistol.'', ''PiÃ±a: Did someone die?'', ''PiÃ±a: Are they loved ones?'', ''Linthe Van Rooij: No'', ''Linthe Van Rooij: Not relevant'', ''PiÃ±a: Were they playing checkers together?'', ''Linthe Van Rooij: Yes'', ''PiÃ±a: Can the winner do something with the pistol?'', ''Linthe Van Rooij: Yes'', ''PiÃ±a: Can the winner kill the loser with the pistol?'', ''Linthe Van Rooij: No'', ''PiÃ±a: Is someone going to shoot with the pistol?'', ''Linthe Van Rooij: Yes'', ''PiÃ±a: Did the woman win?'', ''Linthe Van Rooij: Not relevant'', ''PiÃ±a: Did the man or woman shot another person in the room?'', ''Linthe Van Rooij: No'', ''PiÃ±a: Were there other people in the room?'', ''Linthe Van Rooij: No'', ''PiÃ±a: Did someone shoot in the air?'', ''Linthe Van Rooij: No'', ''PiÃ±a: Did one of the two people kill the other?'', ''Linthe Van Rooij: No'', ''PiÃ±a: I mean did one shot the other?'', ''Linthe Van Rooij: No'', ''PiÃ±a: But there was a shot?'', ''Linthe Van Rooij: No'', ''PiÃ±a: But there was a shot but not on a person?'', ''Linthe Van Rooij: No'', ''Linthe Van Rooij: there was no shot but there is going to be'', ''PiÃ±a: Can the winner use the pistol to make a shot?'', ''Linthe Van Rooij: Yes'', ''PiÃ±a: Will the winner of the game shoot the other player after winning the game of checkers?'', ''Linthe Van Rooij: No'', ''PiÃ±a: Will someone die in the future?'', ''Linthe Van Rooij: Yes'', ''PiÃ±a: One of the two players?'', ''Linthe Van Rooij: Yes'', ''PiÃ±a: The winner?'', ''Linthe Van Rooij: Yes'', ''PiÃ±a: The loser is so sad and mad that they shot the winner?'', ''Linthe Van Rooij: No'', ''PiÃ±a: They are playing the game because someone has to die?'', ''Linthe Van Rooij: No'', ''PiÃ±a: The winner shot themselves?'', ''Linthe Van Rooij: Yes'', ''PiÃ±a: Is the why of this question relevant?'', ''Linthe Van Rooij: Yes'', ''Linthe Van Rooij: Ask about the location'', ''PiÃ±a: Were they in jail?'', ''Linthe Van Rooij: No'', ''PiÃ±a: Were they in a bar?'', ''Linthe Van Rooij: No'', ''PiÃ±a: Were they in an unusual place to play the game?'', ''Linthe Van Rooij: Yes'', ''PiÃ±a: Were they outside?'', ''Linthe Van Rooij: No'', ''PiÃ±a: At someoneâ€™s house?'', ''Linthe Van Rooij: No'', ''PiÃ±a: Were they going from a place to another place while playing the game?'', ''Linthe Van Rooij: No'', ''PiÃ±a: This is so hard, can you give me another hint'', ''Linthe Van Rooij: They are in some kind of vehicle'', ''PiÃ±a: In a plane?'', ''Linthe Van Rooij: No'', ''PiÃ±a: In a train?'', ''Linthe Van Rooij: No'', ''PiÃ±a: In a bus?'', ''Linthe Van Rooij: No'', ''PiÃ±a: In a car?'', ''Linthe Van Rooij: No'', ''PiÃ±a: In a boat?'', ''Linthe Van Rooij: Yes'', ''PiÃ±a: Did someone shoot a hole in the boat?'', ''Linthe Van Rooij: No'', ''Linthe Van Rooij: The type of boat is relevant to the story'', ''PiÃ±a: Okay summary, two players are on a boat playing checkers, the winner fires a shot with the pistol in the future on theirselves'', ''Linthe Van Rooij: Yes'', ''PiÃ±a: Is it a cruise?'', ''Linthe Van Rooij: No'', ''PiÃ±a: A big boat?'', ''Linthe Van Rooij: No'', ''PiÃ±a: Or a small boat for just the two of them?'', ''Linthe Van Rooij: Yes'', ''PiÃ±a: A rowing boat?'', ''Linthe Van Rooij: No'', ''PiÃ±a: A rubber boat?'', ''Linthe Van Rooij: No'', ''PiÃ±a: The boat was big enough for a table right?'', ''Linthe Van Rooij: Yes'', ''PiÃ±a: An expensive boat?'', ''Linthe Van Rooij: Yes'', ''PiÃ±a: A yacht?'', ''Linthe Van Rooij: No'', ''PiÃ±a: A sailing boat?'', ''Linthe Van Rooij: No'', ''PiÃ±a: I just named all boat types haha, do I have to go in a different direction?'', ''Linthe Van Rooij: Yes, under water'', ''PiÃ±a: A submarine?'', ''Linthe Van Rooij: Yes'', ''PiÃ±a: Is why they were there relevant?'', ''Linthe Van Rooij: No'', ''PiÃ±a: What am I missing?'', ''Linthe Van Rooij: Why would you want to shoot yourself in a submarine?'', ''PiÃ±a: Because the person was getting crazy?'', ''Linthe Van Rooij: No'', ''PiÃ±a: Was the submarine still able to go up after the shot?'', ''Linthe Van Rooij: No'', ''PiÃ±a: So the person killed everyone in the submarine by the shot?'', ''Linthe Van Rooij: No'', ''PiÃ±a: Was the shot on themselves intentional?'', ''Linthe Van Rooij: Yes'', ''PiÃ±a: Does the why they wanted to shot themselves have to do with the game?'', ''Linthe Van Rooij: No'', ''PiÃ±a: Did it have something to do with the other person?'', ''Linthe Van Rooij: No'', ''PiÃ±a: Give me more hints please haha I am done with this game'', "Linthe Van Rooij: The game is a result of that the submarine can''t go up", ''PiÃ±a: Am I even close? I am so bad at this'', ''Linthe Van Rooij: Your almost there!!!'', ''PiÃ±a: The person who shot themselves was the driver?'', ''Linthe Van Rooij: Not relevant'', ''PiÃ±a: They were forever stuck there?'', ''Linthe Van Rooij: Yes'', ''PiÃ±a: They already knew that they were stuck before playing the game?'', ''Linthe Van Rooij: Yes'', ''PiÃ±a: There was just one bullet?'', ''Linthe Van Rooij: Yes'', ''PiÃ±a: Ok I think I know now. The people were stuck in a submarine, and they wanted to kill themselves but there was just one bullet in the pistol, they played a game of checkers and the winner could use the bullet to use it to shoot themselves'']'

# Column I (Questions asked) for new rows 16-19
$ws.Cells.Item(16, 9).Value2 = '[''Did she die in a fire or die from the smoke of the fire'', ''Was she already burried?'', ''Was she already presumed to be dead?'', ''Did she fall in a coffin?'', ''But she is in a coffin?'', ''She lit the lighter to see that she was buried?'', ''Was she mistaken for someone else?'', ''Was she in the coffin by accident?'', ''Was she looking for someone in the coffin?'', ''Did people already think that she was dead?'', ''Did she fall asleep in the coffin?'', ''Was there anyone else involved in her death?'', ''Did that person try to kill/murder her or accidentally kill her?'', ''Was that person alive?'', ''Did she suffocate?'', ''Did she dig the hole that she was in herself?'', ''Hint'', ''Was the other person a gravedigger or involved at the graveyard?'', ''Was there anyone else involved'', ''Was the gravedigger digging her grave?'', ''Was the gravedigger already dead when she met him?'', ''Did they fight and she killed him?'', ''Did she know him?'', ''Was she related to him?'', ''Were they partners'', ''To set something on fire or to signal something?'', ''To light a cigarette or to suffocate herself?'', ''Was the other person with her inside or the coffin?'', ''Did she light it to kill him or to do anything related to him?'', ''Did she take a shortcut home through the graveyard?'', ''Did the guy fall in the coffin?'', ''Was she at a funeral?'', ''At a social event?'', ''Was she visiting someone at the graveyard or was she running?'', ''Was she walking'', ''Was she scared of the gravedigger'', ''Was she at home'', ''Anything sport related or food related'', ''Was she drunk or on drugs'', ''Was she shopping or in town'', ''Was it night'', ''Was she driving around'', ''Hint'', ''Did she just rob the house of the gravedigger'', ''Did she rob someone who died'', ''Did she just rob someone'', ''Did she just murder someone?'', ''Was she escaping prison'', ''And she was hiding in the grave to escape the police?'', ''And the gravedigger was an accomplice?'', ''So they tried to hide from the police by hiding in the grave but they got locked in'', ''She tried to hide from the police in the grave that the gravedigger made for her, but he was also dead in the coffin so they were trapped, she lit the lighter to confirm that he was dead'', ''She tried to escape the prison through the coffin of someone else, and the gravedigger was supposed to get her out, but the coffin was for the gravedigger thus he could not get her out, she used the lighter to see the face of the person that was in the coffin with her'']'
$ws.Cells.Item(17, 9).Value2 = '[''is one of them going to die after she leaves?'', ''do they have an important conversation after he received the coke?'', ''does the men also has a weapon?'', ''will he take the knife of the women?'', ''is she leaving with the knife in her hand?'', ''did the men ask the women for a knife?'', ''is one the men or the women getting injured?'', ''are there more people involved except the men & women?'', ''is it important to know why the women decides to leave?'', ''so HE leaves? (not she)'', ''did the  men pay for his coke?'', ''is the relation between the men & women important?'', ''did the women aimed to threaten him with the knife?'', ''did the women ask him to do something with the knife because she did not want to do something?'', ''did the men insult the women?'', ''did the men do something which can be seen considerd as ethically wrong?'', ''did he order a coke to celebrate something?'', ''did he order something because had some medical emergency?'', ''did he order something for a negative reason?'', ''was this negative reason related to the women?'', ''was this related to his work?'', ''to his family live?'', ''is the fact that he order coke (and not another drink) important to understand the negative event?'', ''did he revelaed a secret?'', ''did the men receive his coke before she grabbed the knife?'', ''did the men threatend the women as well?'', ''so did the women refused to give him the coke?'', ''did the negative event happended shortly before he entered the cafe?'', ''did they got into a fight?'', ''did the men drink anything before he order the coke?'', ''were there any other elements or objects involved, caused the negative event?'', ''is the coke dangerous for him?'', ''is the negative event over once he leaves the cafe?'', ''is the women trying to help the men by showing the knife?'', ''does the men wants to kill himself?'', ''does the women wants to kill or injure him?'', ''does the women wants to break or open anything related to the men (rope, hair, a bag ect)?'', ''to protect herself?'', ''to protect and help him?'', ''is the women trying to help the men by scaring him with the knife?'', ''Did the men had a hiccup, and therefore asked for a coke. However, this did not work. Therefore the women tried to scared him with a knife. As a result, the hiccup went away and the men thanked her for her help.'']'
$ws.Cells.Item(18, 9).Value2 = '[''is the vehicle still on land when he climbs out the vehicle?'', ''is the vehicle close to water as well?'', ''were there any other people in the vehicle as well?'', ''was this vehicle an airplane?'', ''is this vehicle used for public transport?'', ''is it relevant to know the type of vehicle?'', ''is the vehicle made to be on land?'', ''is it a tractor?'', ''did he drown himself on purpose?'', ''was he a farmer?'', ''did something negative happened to his farming practices?'', ''did something negative happened at all?'', ''does this has anything to do with the tractor?'', ''was the tractor broken?'', ''was the land, in which the tractor was located at that moment boggy?'', ''did the reason that he climbed out of the tractor had anything to do with the tractor?'', ''was the water in which he drowned himself close to the tractor?'', ''were there any other objects or people involved in the negative event?'', ''were there any people in the tractor?'', ''outside the tractor?'', ''did he wanted to go to the other people?'', ''did the drowning happened right after he climbed out of the tractor?'', ''did the other people influenced him in stepping out of the tractor?'', ''was it dangerous for him to stay in the vehicle?'', ''were the other people also with tractors? n'', ''is it important to know how he exactly climbed out of the tractor?'', ''did the other people wanted him to drown himself?'', ''did the other people wanted to save him?'', ''did he know the other people before?'', ''was it an accident that he drowned?'', ''did the other people saw something which the men could not see (when he was still in the tractor)?'', ''were the other people also farmers?'', ''were the other people angry at him?'', ''did the know eachother for a long time?'', ''were they employees?'', ''family?'', ''were they close to the tractor?'', ''was there a fight or any other negative event between the men and his family?'', ''were there any childern involved?'', ''were the childeren outside the tractor?'', ''did he climbed out the tractor for his childern?'', ''was his kid in a dangerous situation?'', ''did he aimed to save his kid by climbing out of the tractor?'', ''did he cause the dangerous situation?'', ''did he accidently killed his kid with his tractor?'', ''so he accidently he run over the kid with his tractor because he did not saw him because they were playing the field so he drowned himself somewhere afterwards because it made him so depressed'']'
$ws.Cells.Item(19, 9).Value2 = '[''Did someone die?'', ''Are they loved ones?'', ''Were they playing checkers together?'', ''Can the winner do something with the pistol?'', ''Can the winner kill the loser with the pistol?'', ''Is someone going to shoot with the pistol?'', ''Did the woman win?'', ''Did the man or woman shot another person in the room?'', ''Were there other people in the room?'', ''Did someone shoot in the air?'', ''Did one of the two people kill the other?'', ''I mean did one shot the other?'', ''But there was a shot?'', ''But there was a shot but not on a person?'', ''Can the winner use the pistol to make a shot?'', ''Will the winner of the game shoot the other player after winning the game of checkers?'', ''Will someone die in the future?'', ''One of the two players?'', ''The winner?'', ''The loser is so sad and mad that they shot the winner?'', ''They are playing the game because someone has to die?'', ''The winner shot themselves?'', ''Is the why of this question relevant?'', ''Were they in jail?'', ''Were they in a bar?'', ''Were they in an unusual place to play the game?'', ''Were they outside?'', ''At someoneâ€™s house?'', ''Were they going from a place to another place while playing the game?'', ''This is so hard, can you give me another hint'', ''In a plane?'', ''In a train?'', ''In a bus?'', ''In a car?'', ''In a boat?'', ''Did someone shoot a hole in the boat?'', ''Okay summary, two players are on a boat playing checkers, the winner fires a shot with the pistol in the future on theirselves'', ''Is it a cruise?'', ''A big boat?'', ''Or a small boat for just the two of them?'', ''A rowing boat?'', ''A rubber boat?'', ''The boat was big enough for a table right?'', ''An expensive boat?'', ''A yacht?'', ''A sailing boat?'', ''I just named all boat types haha, do I have to go in a different direction?'', ''A submarine?'', ''Is why they were there relevant?'', ''What am I missing?'', ''Because the person was getting crazy?'', ''Was the submarine still able to go up after the shot?'', ''So the person killed everyone in the submarine by the shot?'', ''Was the shot on themselves intentional?'', ''Does the why they wanted to shot themselves have to do with the game?'', ''Did it have something to do with the other person?'', ''Give me more hints please haha I am done with this game'', ''Am I even close? I am so bad at this'', ''The person who shot themselves was the driver?'', ''They were forever stuck there?'', ''They already knew that they were stuck before playing the game?'', ''There was just one bullet?'', ''Ok I think I know now. The people were stuck in a submarine, and they wanted to kill themselves but there was just one bullet in the pistol, they played a game of checkers and the winner could use the bullet to use it to shoot themselves'']'

# Columns J-K (numeric) for new rows 16-19
$ws.Cells.Item(16, 10).Value2 = 10
$ws.Cells.Item(16, 11).Value2 = 40
$ws.Cells.Item(17, 10).Value2 = 7
$ws.Cells.Item(17, 11).Value2 = 33
$ws.Cells.Item(18, 10).Value2 = 22
$ws.Cells.Item(18, 11).Value2 = 23
$ws.Cells.Item(19, 10).Value2 = 21
$ws.Cells.Item(19, 11).Value2 = 38

# Column L (Hints) for new rows 16-19
$ws.Cells.Item(16, 12).Value2 = '[''A woman dies after lighting a lighter, and doesnâ€™t even need a funeral.'', ''Ask about the occupation of the other person involved'', ''Why did she lit the lighter inside the coffin?'', ''Ask about where she was before'', ''Not relevant'', ''She was a criminal/convict'', ''She expected another person to be in the grave next to her'', ''Allrightyy'']'
$ws.Cells.Item(17, 12).Value2 = '[''A man goes into a cafe and orders a glass of coke. The woman behind the bar grabs a knife and pointed it towards the man. He thanks her and leaves.'', ''False assumption woman leaves'', ''Not relevant'', ''Ask about why he orders something at the bar'', ''No, it is still happening'', ''The drink could serve as a remedy for the event'', ''She only wants to scare him'', ''Ask these questions seperately'', ''Yes!!!!!!! You got it'']'
$ws.Cells.Item(18, 12).Value2 = '[''A man climbs out of a large vehicle and drowned himself'', ''Not relevant'', ''Not relevant'', ''Not relevant'', ''Ask about the relationship with the people involved'', ''Not relevant'']'
$ws.Cells.Item(19, 12).Value2 = '[''A man and a woman are seated at a table. On the table between them is a game of checkers and a pistol.'', ''Not relevant'', ''Not relevant'', ''there was no shot but there is going to be'', ''Ask about the location'', ''They are in some kind of vehicle'', ''The type of boat is relevant to the story'', ''Yes, under water'', ''Why would you want to shoot yourself in a submarine?'', "The game is a result of that the submarine can''t go up", ''Your almost there!!!'', ''Not relevant'']'

# Update selection to K13 (matches final sheetView state)
$ws.Range("K13").Select()

# Page setup changes
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
